$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 117, pushing existing rows 117-170 down to 118-171
$ws.Rows("117:117").Insert()

# Populate the new row 117 with the new price record
$ws.Range("A117").Value = 1
$ws.Range("B117").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C117").Value = "Arica y Parinacota"
$ws.Range("D117").Value = 45007
$ws.Range("E117").Value = 15
$ws.Range("F117").Value = "Fruta"
$ws.Range("G117").Value = 100106
$ws.Range("H117").Value = "Oleaginosos"
$ws.Range("I117").Value = 100106002
$ws.Range("J117").Value = "Palta"
$ws.Range("K117").Value = "Hass"
$ws.Range("L117").Value = "Primera"
$ws.Range("M117").Value = 400
$ws.Range("N117").Value = 30000
$ws.Range("O117").Value = 31000
$ws.Range("P117").Value = 30500
$ws.Range("Q117").Value = "$/bandeja 10 kilos"
$ws.Range("R117").Value = "Perú"
$ws.Range("S117").Value = 3050
$ws.Range("T117").Value = 10
